# Update NATMI TPM-derived metrics (Sema4a-Plxnb1 LR pair) to reflect new TPM inputs.
# Each assignment below corresponds to one changed <c r="...">/<v> cell from the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 6.051203666666667
$ws.Range("H2").Value = 18.153611
$ws.Range("I2").Value = 0.07608037240065801
$ws.Range("J2").Value = 0.07775008964215516
$ws.Range("M2").Value = 0.600843
$ws.Range("N2").Value = 1.802529
$ws.Range("O2").Value = 0.1216566842860732
$ws.Range("P2").Value = 0.148762828988552
$ws.Range("Q2").Value = 3.635823364691
$ws.Range("R2").Value = 32.722410282219
$ws.Range("S2").Value = 0.009255685845513728
$ws.Range("T2").Value = 0.01156632328928052

# Row 3
$ws.Range("G3").Value = 6.051203666666667
$ws.Range("H3").Value = 18.153611
$ws.Range("I3").Value = 0.07608037240065801
$ws.Range("J3").Value = 0.07775008964215516
$ws.Range("O3").Value = 0.2991130341144489
$ws.Range("P3").Value = 0.3657579639239645
$ws.Range("Q3").Value = 8.939271725996113
$ws.Range("R3").Value = 80.45344553396501
$ws.Range("S3").Value = 0.022756631025318
$ws.Range("T3").Value = 0.02843771448242039

# Row 4
$ws.Range("G4").Value = 6.051203666666667
$ws.Range("H4").Value = 18.153611
$ws.Range("I4").Value = 0.07608037240065801
$ws.Range("J4").Value = 0.07775008964215516
$ws.Range("M4").Value = 0.111967
$ws.Range("N4").Value = 0.335901
$ws.Range("O4").Value = 0.02267070427625646
$ws.Range("P4").Value = 0.02772193014375004
$ws.Range("Q4").Value = 0.6775351209456667
$ws.Range("R4").Value = 6.097816088511
$ws.Range("S4").Value = 0.001724795623922781
$ws.Range("T4").Value = 0.002155382553730129

# Row 5
$ws.Range("G5").Value = 6.051203666666667
$ws.Range("H5").Value = 18.153611
$ws.Range("I5").Value = 0.07608037240065801
$ws.Range("J5").Value = 0.07775008964215516
$ws.Range("M5").Value = 2.6997255
$ws.Range("N5").Value = 5.399451
$ws.Range("O5").Value = 0.5466314042313235
$ws.Range("P5").Value = 0.4456170223863617
$ws.Range("Q5").Value = 16.3365888445935
$ws.Range("R5").Value = 98.01953306756101
$ws.Range("S5").Value = 0.04158792079981372
$ws.Range("T5").Value = 0.03464676343660988

# Row 6
$ws.Range("G6").Value = 6.051203666666667
$ws.Range("H6").Value = 18.153611
$ws.Range("I6").Value = 0.07608037240065801
$ws.Range("J6").Value = 0.07775008964215516
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.04903366666666667
$ws.Range("N6").Value = 0.147101
$ws.Range("O6").Value = 0.009928173091897913
$ws.Range("P6").Value = 0.01214025455737189
$ws.Range("Q6").Value = 0.2967127035234445
$ws.Range("R6").Value = 2.670414331711001
$ws.Range("S6").Value = 0.0007553391060897856
$ws.Range("T6").Value = 0.0009439058801142471

# Row 7
$ws.Range("I7").Value = 0.2215826302097334
$ws.Range("J7").Value = 0.2264456497560767
$ws.Range("M7").Value = 0.600843
$ws.Range("N7").Value = 1.802529
$ws.Range("O7").Value = 0.1216566842860732
$ws.Range("P7").Value = 0.148762828988552
$ws.Range("Q7").Value = 10.589266044645
$ws.Range("R7").Value = 95.30339440180501
$ws.Range("S7").Value = 0.02695700808670324
$ws.Range("T7").Value = 0.03368669546986479

# Row 8
$ws.Range("I8").Value = 0.2215826302097334
$ws.Range("J8").Value = 0.2264456497560767
$ws.Range("O8").Value = 0.2991130341144489
$ws.Range("P8").Value = 0.3657579639239645
$ws.Range("S8").Value = 0.06627825282909329
$ws.Range("T8").Value = 0.0828242997942218

# Row 9
$ws.Range("I9").Value = 0.2215826302097334
$ws.Range("J9").Value = 0.2264456497560767
$ws.Range("M9").Value = 0.111967
$ws.Range("N9").Value = 0.335901
$ws.Range("O9").Value = 0.02267070427625646
$ws.Range("P9").Value = 0.02772193014375004
$ws.Range("Q9").Value = 1.973308087505
$ws.Range("R9").Value = 17.759772787545
$ws.Range("S9").Value = 0.005023434282239955
$ws.Range("T9").Value = 0.006277510483894046

# Row 10
$ws.Range("I10").Value = 0.2215826302097334
$ws.Range("J10").Value = 0.2264456497560767
$ws.Range("M10").Value = 2.6997255
$ws.Range("N10").Value = 5.399451
$ws.Range("O10").Value = 0.5466314042313235
$ws.Range("P10").Value = 0.4456170223863617
$ws.Range("Q10").Value = 47.5800027078825
$ws.Range("R10").Value = 285.480016247295
$ws.Range("S10").Value = 0.1211240243048166
$ws.Range("T10").Value = 0.1009080361766479

# Row 11
$ws.Range("I11").Value = 0.2215826302097334
$ws.Range("J11").Value = 0.2264456497560767
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.04903366666666667
$ws.Range("N11").Value = 0.147101
$ws.Range("O11").Value = 0.009928173091897913
$ws.Range("P11").Value = 0.01214025455737189
$ws.Range("Q11").Value = 0.8641700768383334
$ws.Range("R11").Value = 7.777530691545
$ws.Range("S11").Value = 0.00219991070688024
$ws.Range("T11").Value = 0.002749107831448249

# Row 12
$ws.Range("G12").Value = 33.62840566666667
$ws.Range("H12").Value = 100.885217
$ws.Range("I12").Value = 0.4228021014155913
$ws.Range("J12").Value = 0.432081235260482
$ws.Range("M12").Value = 0.600843
$ws.Range("N12").Value = 1.802529
$ws.Range("O12").Value = 0.1216566842860732
$ws.Range("P12").Value = 0.148762828988552
$ws.Range("Q12").Value = 20.205392145977
$ws.Range("R12").Value = 181.848529313793
$ws.Range("S12").Value = 0.05143670176740489
$ws.Range("T12").Value = 0.06427762691021741

# Row 13
$ws.Range("G13").Value = 33.62840566666667
$ws.Range("H13").Value = 100.885217
$ws.Range("I13").Value = 0.4228021014155913
$ws.Range("J13").Value = 0.432081235260482
$ws.Range("O13").Value = 0.2991130341144489
$ws.Range("P13").Value = 0.3657579639239645
$ws.Range("Q13").Value = 49.67829088653944
$ws.Range("R13").Value = 447.104617978855
$ws.Range("S13").Value = 0.1264656193843824
$ws.Range("T13").Value = 0.1580371528586254

# Row 14
$ws.Range("G14").Value = 33.62840566666667
$ws.Range("H14").Value = 100.885217
$ws.Range("I14").Value = 0.4228021014155913
$ws.Range("J14").Value = 0.432081235260482
$ws.Range("M14").Value = 0.111967
$ws.Range("N14").Value = 0.335901
$ws.Range("O14").Value = 0.02267070427625646
$ws.Range("P14").Value = 0.02772193014375004
$ws.Range("Q14").Value = 3.765271697279666
$ws.Range("R14").Value = 33.887445275517
$ws.Range("S14").Value = 0.009585221408572661
$ws.Range("T14").Value = 0.01197812582031631

# Row 15
$ws.Range("G15").Value = 33.62840566666667
$ws.Range("H15").Value = 100.885217
$ws.Range("I15").Value = 0.4228021014155913
$ws.Range("J15").Value = 0.432081235260482
$ws.Range("M15").Value = 2.6997255
$ws.Range("N15").Value = 5.399451
$ws.Range("O15").Value = 0.5466314042313235
$ws.Range("P15").Value = 0.4456170223863617
$ws.Range("Q15").Value = 90.78746430264449
$ws.Range("R15").Value = 544.724785815867
$ws.Range("S15").Value = 0.2311169064087591
$ws.Range("T15").Value = 0.1925427534857971

# Row 16
$ws.Range("G16").Value = 33.62840566666667
$ws.Range("H16").Value = 100.885217
$ws.Range("I16").Value = 0.4228021014155913
$ws.Range("J16").Value = 0.432081235260482
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.04903366666666667
$ws.Range("N16").Value = 0.147101
$ws.Range("O16").Value = 0.009928173091897913
$ws.Range("P16").Value = 0.01214025455737189
$ws.Range("Q16").Value = 1.648924033990778
$ws.Range("R16").Value = 14.840316305917
$ws.Range("S16").Value = 0.004197652446472166
$ws.Range("T16").Value = 0.005245576185525943

# Row 17
$ws.Range("G17").Value = 5.124275
$ws.Range("H17").Value = 10.24855
$ws.Range("I17").Value = 0.06442631445887793
$ws.Range("J17").Value = 0.04389350863594627
$ws.Range("M17").Value = 0.600843
$ws.Range("N17").Value = 1.802529
$ws.Range("O17").Value = 0.1216566842860732
$ws.Range("P17").Value = 0.148762828988552
$ws.Range("Q17").Value = 3.078884763825
$ws.Range("R17").Value = 18.47330858295
$ws.Range("S17").Value = 0.007837891797838985
$ws.Range("T17").Value = 0.006529722518916806

# Row 18
$ws.Range("G18").Value = 5.124275
$ws.Range("H18").Value = 10.24855
$ws.Range("I18").Value = 0.06442631445887793
$ws.Range("J18").Value = 0.04389350863594627
$ws.Range("O18").Value = 0.2991130341144489
$ws.Range("P18").Value = 0.3657579639239645
$ws.Range("Q18").Value = 7.569946269708334
$ws.Range("R18").Value = 45.41967761825
$ws.Range("S18").Value = 0.01927075039460657
$ws.Range("T18").Value = 0.01605440034816266

# Row 19
$ws.Range("G19").Value = 5.124275
$ws.Range("H19").Value = 10.24855
$ws.Range("I19").Value = 0.06442631445887793
$ws.Range("J19").Value = 0.04389350863594627
$ws.Range("M19").Value = 0.111967
$ws.Range("N19").Value = 0.335901
$ws.Range("O19").Value = 0.02267070427625646
$ws.Range("P19").Value = 0.02772193014375004
$ws.Range("Q19").Value = 0.5737496989249999
$ws.Range("R19").Value = 3.44249819355
$ws.Range("S19").Value = 0.001460589922706327
$ws.Range("T19").Value = 0.001216812780169791

# Row 20
$ws.Range("G20").Value = 5.124275
$ws.Range("H20").Value = 10.24855
$ws.Range("I20").Value = 0.06442631445887793
$ws.Range("J20").Value = 0.04389350863594627
$ws.Range("M20").Value = 2.6997255
$ws.Range("N20").Value = 5.399451
$ws.Range("O20").Value = 0.5466314042313235
$ws.Range("P20").Value = 0.4456170223863617
$ws.Range("Q20").Value = 13.8341358865125
$ws.Range("R20").Value = 55.33654354605
$ws.Range("S20").Value = 0.03521744674210527
$ws.Range("T20").Value = 0.01955969462044043

# Row 21
$ws.Range("G21").Value = 5.124275
$ws.Range("H21").Value = 10.24855
$ws.Range("I21").Value = 0.06442631445887793
$ws.Range("J21").Value = 0.04389350863594627
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 0.6666666666666666
$ws.Range("M21").Value = 0.04903366666666667
$ws.Range("N21").Value = 0.147101
$ws.Range("O21").Value = 0.009928173091897913
$ws.Range("P21").Value = 0.01214025455737189
$ws.Range("Q21").Value = 0.2512619922583333
$ws.Range("R21").Value = 1.50757195355
$ws.Range("S21").Value = 0.0006396356016207854
$ws.Range("T21").Value = 0.000532878368256589

# Row 22
$ws.Range("G22").Value = 17.10908866666667
$ws.Range("H22").Value = 51.327266
$ws.Range("I22").Value = 0.2151085815151395
$ws.Range("J22").Value = 0.2198295167053399
$ws.Range("M22").Value = 0.600843
$ws.Range("N22").Value = 1.802529
$ws.Range("O22").Value = 0.1216566842860732
$ws.Range("P22").Value = 0.148762828988552
$ws.Range("Q22").Value = 10.279876161746
$ws.Range("R22").Value = 92.518885455714
$ws.Range("S22").Value = 0.02616939678861236
$ws.Range("T22").Value = 0.03270246080027252

# Row 23
$ws.Range("G23").Value = 17.10908866666667
$ws.Range("H23").Value = 51.327266
$ws.Range("I23").Value = 0.2151085815151395
$ws.Range("J23").Value = 0.2198295167053399
$ws.Range("O23").Value = 0.2991130341144489
$ws.Range("P23").Value = 0.3657579639239645
$ws.Range("Q23").Value = 25.27477192975445
$ws.Range("R23").Value = 227.47294736779
$ws.Range("S23").Value = 0.06434178048104863
$ws.Range("T23").Value = 0.08040439644053426

# Row 24
$ws.Range("G24").Value = 17.10908866666667
$ws.Range("H24").Value = 51.327266
$ws.Range("I24").Value = 0.2151085815151395
$ws.Range("J24").Value = 0.2198295167053399
$ws.Range("M24").Value = 0.111967
$ws.Range("N24").Value = 0.335901
$ws.Range("O24").Value = 0.02267070427625646
$ws.Range("P24").Value = 0.02772193014375004
$ws.Range("Q24").Value = 1.915653330740667
$ws.Range("R24").Value = 17.240879976666
$ws.Range("S24").Value = 0.004876663038814733
$ws.Range("T24").Value = 0.006094098505639764

# Row 25
$ws.Range("G25").Value = 17.10908866666667
$ws.Range("H25").Value = 51.327266
$ws.Range("I25").Value = 0.2151085815151395
$ws.Range("J25").Value = 0.2198295167053399
$ws.Range("M25").Value = 2.6997255
$ws.Range("N25").Value = 5.399451
$ws.Range("O25").Value = 0.5466314042313235
$ws.Range("P25").Value = 0.4456170223863617
$ws.Range("Q25").Value = 46.18984295516101
$ws.Range("R25").Value = 277.139057730966
$ws.Range("S25").Value = 0.1175851059758288
$ws.Range("T25").Value = 0.09795977466686653

# Row 26
$ws.Range("G26").Value = 17.10908866666667
$ws.Range("H26").Value = 51.327266
$ws.Range("I26").Value = 0.2151085815151395
$ws.Range("J26").Value = 0.2198295167053399
$ws.Range("K26").Value = 2
$ws.Range("L26").Value = 0.6666666666666666
$ws.Range("M26").Value = 0.04903366666666667
$ws.Range("N26").Value = 0.147101
$ws.Range("O26").Value = 0.009928173091897913
$ws.Range("P26").Value = 0.01214025455737189
$ws.Range("Q26").Value = 0.8389213506517779
$ws.Range("R26").Value = 7.550292155866001
$ws.Range("S26").Value = 0.002135635230834937
$ws.Range("T26").Value = 0.002668786292026863
